$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Header row (row 1): extend the table with the standardized metadata
# columns already used on the other sheets (land/building/car/deposit/
# stock/fund). Copy the existing bold/bordered header style (s=1, taken
# from E1) across the new header cells before writing their text so the
# formatting matches the rest of row 1.
$ws.Range("E1").Copy($ws.Range("F1:K1"))
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Data rows (2-5): the old E column held a raw acquisition amount;
# it is replaced by the property_category marker ("insurance"), and the
# row is extended with category/date/legislator_name/legislator_id/
# source_file/index, matching the other sheets.
$rows = @(
    @{ Row = 2; Index = 123 },
    @{ Row = 3; Index = 124 },
    @{ Row = 4; Index = 125 },
    @{ Row = 5; Index = 126 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("E$row").Value = "insurance"
    $ws.Range("F$row").Value = "normal"

    # Writing a date-shaped string via .Value auto-converts it to a date
    # serial; force text first so it is stored the same way as the other
    # sheets (a literal "2012-04-26" string), then drop back to the plain
    # "Normal" style so no stray number-format sticks to the cell.
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = "2012-04-26"
    $ws.Range("G$row").Style = "Normal"

    $ws.Range("H$row").Value = "王進士"
    $ws.Range("I$row").Value = 1701
    $ws.Range("J$row").Value = "tmp4a4e1"
    $ws.Range("K$row").Value = $r.Index
}
